$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price/Volume columns so numeric-looking values
# (e.g. "1.001", "304.20") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.475.50'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.631.93'
$ws.Range("E3").Value = '  -0.40%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '0.9997'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("D6").Value = '304.20'
$ws.Range("E6").Value = '  -1.34%  '

$ws.Range("D7").Value = '0.3771'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").Value = '51.90'
$ws.Range("E8").Value = '  -2.06%  '

$ws.Range("D9").Value = '0.3642'
$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.08155'
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '1.231'
$ws.Range("E11").Value = '  -3.13%  '

$ws.Range("D12").Value = '0.9993'
$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("E13").Value = '  -1.71%  '

$ws.Range("D14").Value = '6.571'
$ws.Range("E14").Value = '  -1.55%  '

$ws.Range("D15").Value = '0.00001251'
$ws.Range("E15").Value = '  -2.22%  '

$ws.Range("D16").Value = '7.263'
$ws.Range("E16").Value = '  -2.68%  '

$ws.Range("D17").Value = '1.629.07'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '93.98'
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").Value = '0.06954'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").Value = '17.89'
$ws.Range("E20").Value = '  -2.41%  '

$ws.Range("D21").Value = '6.434'
$ws.Range("E21").Value = '  -2.26%  '

$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '12.76'
$ws.Range("E23").Value = '  -0.97%  '

$ws.Range("B24").Value = 'WrappedBTC'
$ws.Range("C24").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D24").Value = '23.475.51'
$ws.Range("E24").Value = '  -0.14%  '

$ws.Range("D25").Value = '3.217'
$ws.Range("E25").Value = '  +4.41%  '

$ws.Range("D26").Value = '2.447'
$ws.Range("E26").Value = '  +1.27%  '

$ws.Range("D27").Value = '21.28'
$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("D28").Value = '150.45'
$ws.Range("E28").Value = '  -0.29%  '

$ws.Range("D29").Value = '5.280'
$ws.Range("E29").Value = '  -0.67%  '

$ws.Range("D30").Value = '135.15'
$ws.Range("E30").Value = '  -0.83%  '

$ws.Range("B31").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C31").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D31").Value = '1.805.94'
$ws.Range("E31").Value = '  -0.65%  '

$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.314'
$ws.Range("E32").Value = '  -3.93%  '

$ws.Range("D33").Value = '6.844'
$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '0.9818'
$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = '10.96'
$ws.Range("E35").Value = '  +5.17%  '

$ws.Range("D36").Value = '0.02792'
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").Value = '0.2543'
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.07218'
$ws.Range("E38").Value = '  -2.77%  '

$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").Value = '0.08832'
$ws.Range("E39").Value = '  -0.63%  '

$ws.Range("D40").Value = '6.103'
$ws.Range("E40").Value = '  -0.98%  '

$ws.Range("D41").Value = '0.7090'
$ws.Range("E41").Value = '  -0.62%  '

$ws.Range("D42").Value = '1.358'
$ws.Range("E42").Value = '  -2.56%  '

$ws.Range("D43").Value = '16.35'
$ws.Range("E43").Value = '  -0.23%  '

$ws.Range("D44").Value = '12.37'
$ws.Range("E44").Value = '  -1.53%  '

$ws.Range("D45").Value = '0.6541'
$ws.Range("E45").Value = '  -0.74%  '

$ws.Range("D46").Value = '2.339'
$ws.Range("E46").Value = '  -0.89%  '

$ws.Range("D47").Value = '0.9986'
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").Value = '3.998'
$ws.Range("E48").Value = '  -1.14%  '

$ws.Range("D49").Value = '0.08015'
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("E50").Value = '  -0.62%  '

$ws.Range("D51").Value = '125.75'
$ws.Range("E51").Value = '  -3.52%  '

# Restore original (default) style now that values are safely stored as text.
$ws.Range("D2:E51").Style = "Normal"
